$d = $word.ActiveDocument

# 2) "Но они я – да." -> "Но они я - да." (em dash -> hyphen)
$d.Content.Find.Execute("Но они я – да.", $true, $false, $false, $false, $false, $true, 1, $false, "Но они я - да.", 2)

# 3) "проектом «Хранители». " -> "проектом “Хранители”." (curly quotes, drop extra space)
$d.Content.Find.Execute("заинтересовался проектом «Хранители». ", $true, $false, $false, $false, $false, $true, 1, $false, "заинтересовался проектом “Хранители”.", 2)

# 4) "персонаж, так Тони Старк." -> "персонаж, как Тони Старк."
$d.Content.Find.Execute("персонаж, так Тони Старк.", $true, $false, $false, $false, $false, $true, 1, $false, "персонаж, как Тони Старк.", 2)

# 5) "Привет, Мэри Джейн, – мой ... откровенным. – Восхитительно выглядишь." -> hyphen variants
$d.Content.Find.Execute("Привет, Мэри Джейн, – мой судорожно-восхищенный вздох получился слишком откровенным. – Восхитительно выглядишь.", $true, $false, $false, $false, $false, $true, 1, $false, "Привет, Мэри Джейн,  мой судорожно-восхищенный вздох получился слишком откровенным. - Восхитительно выглядишь.", 2)
